$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Refresh the "datetimeFigureOut" date placeholder on the slide master and
#    on every custom (slide) layout: 1/24/24 -> 1/30/24
# ---------------------------------------------------------------------------
function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.Name -like "Date Placeholder*" -and $sh.HasTextFrame) {
            $tr = $sh.TextFrame.TextRange
            if ($tr.Text -eq "1/24/24") {
                $tr.Text = "1/30/24"
            }
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    Update-DatePlaceholder $layout.Shapes
}

# ---------------------------------------------------------------------------
# 2) Reposition the eight "excluded participants" callout textboxes that sit
#    under the workflow diagram boxes (counts on excluded adults/children) so
#    they line up under their boxes again. Only the vertical (Top) position
#    changes; identified by shape Id (name is duplicated: "TextBox 50").
#
#    Shape.Top is a single-precision (points) property, so the literals below
#    are the exact point values whose float32 round-trip reproduces the
#    target EMU offset from the source file (target_emu = 914400/72*pt).
# ---------------------------------------------------------------------------
$newTop = @{
    5  = 215.70709228515625   # 2739480 EMU
    7  = 192.6844940185547    # 2447093 EMU
    8  = 169.6618194580078    # 2154705 EMU
    9  = 192.6844940185547    # 2447093 EMU
    10 = 169.6618194580078    # 2154705 EMU
    11 = 192.6844940185547    # 2447093 EMU
    12 = 169.6618194580078    # 2154705 EMU
    13 = 192.6844940185547    # 2447093 EMU
}

$slide = $p.Slides.Item(1)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $sh = $slide.Shapes.Item($i)
    if ($newTop.ContainsKey($sh.Id)) {
        $sh.Top = $newTop[$sh.Id]
    }
}
